$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.433.21'
$ws.Range('E2').Value = '  +2.65%  '
$ws.Range('D3').Value = '3.265.15'
$ws.Range('E3').Value = '  -0.45%  '
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').Value = '575.76'
$ws.Range('E5').Value = '  -0.35%  '
$ws.Range('E6').Value = '  -3.33%  '
$ws.Range('E7').Value = '  +0.11%  '
$ws.Range('E8').Value = '  +2.03%  '
$ws.Range('D9').Value = '3.260.19'
$ws.Range('E9').Value = '  -0.46%  '
$ws.Range('D10').Value = '0.175'
$ws.Range('E10').Value = '  -0.26%  '
$ws.Range('D11').Value = '0.571'
$ws.Range('E11').Value = '  -0.04%  '
$ws.Range('D12').Value = '45.28'
$ws.Range('E12').Value = '  -2.13%  '
$ws.Range('E13').Value = '  +1.95%  '
$ws.Range('D14').Value = '679.97'
$ws.Range('E14').Value = '  +9.79%  '
$ws.Range('D15').Value = '3.787.69'
$ws.Range('E15').Value = '  -0.45%  '
$ws.Range('D16').Value = '8.31'
$ws.Range('E16').Value = '  -1.34%  '
$ws.Range('D17').Value = '67.504.68'
$ws.Range('E17').Value = '  +2.87%  '
$ws.Range('E18').Value = '  +1.26%  '
$ws.Range('D19').Value = '3.263.79'
$ws.Range('E19').Value = '  -0.54%  '
$ws.Range('D20').Value = '17.31'
$ws.Range('E20').Value = '  -2.80%  '
$ws.Range('D21').Value = '10.66'
$ws.Range('E21').Value = '  -2.36%  '
$ws.Range('D22').Value = '0.887'
$ws.Range('E22').Value = '  -0.27%  '
$ws.Range('D23').Value = '16.98'
$ws.Range('E23').Value = '  -5.87%  '
$ws.Range('D24').Value = '5.13'
$ws.Range('E24').Value = '  +3.16%  '
$ws.Range('D25').Value = '98.19'
$ws.Range('E25').Value = '  -2.73%  '
$ws.Range('E26').Value = '  -2.28%  '
$ws.Range('E27').Value = '  +0.50%  '
$ws.Range('E28').Value = '  -2.11%  '
$ws.Range('D29').Value = '32.51'
$ws.Range('E29').Value = '  +5.29%  '
$ws.Range('D30').Value = '8.38'
$ws.Range('E30').Value = '  -0.54%  '
$ws.Range('E31').Value = '  +3.01%  '
$ws.Range('D32').Value = '580.30'
$ws.Range('E32').Value = '  +6.08%  '
$ws.Range('D33').Value = '3.867.10'
$ws.Range('E33').Value = '  +2.01%  '
$ws.Range('D34').Value = '10.79'
$ws.Range('E34').Value = '  -0.54%  '
$ws.Range('E35').Value = '  -0.21%  '
$ws.Range('D36').Value = '0.998'
$ws.Range('E36').Value = '  -0.14%  '
$ws.Range('D37').Value = '3.35'
$ws.Range('E37').Value = '  -10.14%  '
$ws.Range('E38').Value = '  -1.63%  '
$ws.Range('E39').Value = '  +1.16%  '
$ws.Range('E40').Value = '  +2.33%  '
$ws.Range('D41').Value = '3.43'
$ws.Range('E41').Value = '  +1.36%  '
$ws.Range('E42').Value = '  +1.06%  '
$ws.Range('D43').Value = '32.08'
$ws.Range('E43').Value = '  -1.27%  '
$ws.Range('D44').Value = '0.0₃0667'
$ws.Range('E44').Value = '  -1.49%  '
$ws.Range('E45').Value = '  -0.44%  '
$ws.Range('E46').Value = '  +0.98%  '
$ws.Range('D47').Value = '0.128'
$ws.Range('E47').Value = '  +1.14%  '
$ws.Range('E48').Value = '  +0.48%  '
$ws.Range('E49').Value = '  +9.21%  '
$ws.Range('E50').Value = '  -0.47%  '
$ws.Range('D51').Value = '129.30'
$ws.Range('E51').Value = '  +0.68%  '
